$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 82
$ws1.Range("F4").Value = 251
$ws1.Range("F6").Value = 10055
$ws1.Range("F8").Value = 907
$ws1.Range("F9").Value = 1250
$ws1.Range("F10").Value = 5544
$ws1.Range("D15").Value = "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"
$ws1.Range("F15").Value = 3077
$ws1.Range("F18").Value = 593
$ws1.Range("F20").Value = 14
$ws1.Range("F22").Value = 17
$ws1.Range("F23").Value = 1523

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 82
$ws4.Range("F5").Value = 251
$ws4.Range("F7").Value = 10055
$ws4.Range("F9").Value = 907
$ws4.Range("F10").Value = 1250
$ws4.Range("F11").Value = 5544
$ws4.Range("D16").Value = "凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"
$ws4.Range("F16").Value = 3077
$ws4.Range("F19").Value = 593
$ws4.Range("F21").Value = 14
$ws4.Range("F23").Value = 17
$ws4.Range("F24").Value = 1523
